$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2030773565271132
$ws.Range("C2").Value = 0.3586383285908433
$ws.Range("D2").Value = 0.2920997974838273
$ws.Range("E2").Value = -0.04706463449054476
$ws.Range("F2").Value = -0.02521342967184042
$ws.Range("G2").Value = -0.01456839544959351
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.4310742903707728
$ws.Range("J2").Value = -0.008345833758740517

$ws.Range("B3").Value = 0.3060149997554619
$ws.Range("C3").Value = 0.6045913867399139
$ws.Range("D3").Value = 0.4395497136685166
$ws.Range("E3").Value = 0.606274581506411
$ws.Range("F3").Value = 0.4632617221138344
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.4350634056425048
$ws.Range("J3").Value = 0.4862065791702738

$ws.Range("B4").Value = 0.2035648004861985
$ws.Range("C4").Value = 0.03759293020688034
$ws.Range("D4").Value = 0.07975460823787993
$ws.Range("E4").Value = 0.3958446515967743
$ws.Range("F4").Value = 0.2319513376959622
$ws.Range("G4").Value = 0.1086963201998874
$ws.Range("H4").Value = 0.08424470825742002
$ws.Range("I4").Value = 0.1120178051926603
$ws.Range("J4").Value = 0.1848811027798208

$ws.Range("B5").Value = 0.8491860042760604
$ws.Range("C5").Value = 0.7103531195001243
$ws.Range("D5").Value = 0.2379363358182544
$ws.Range("E5").Value = 0.6147161768971295
$ws.Range("F5").Value = 0.7780876219875478
$ws.Range("G5").Value = 0.7439439017723934
$ws.Range("H5").Value = 0.05129425985305078
$ws.Range("I5").Value = 0.7780876219875478
$ws.Range("J5").Value = 0.5843133674214755

$ws.Range("B6").Value = 0.2189349112426036
$ws.Range("C6").Value = 0.08284023668639054
$ws.Range("D6").Value = 0.3155576529809658
$ws.Range("E6").Value = 0.2534435261707989
$ws.Range("F6").Value = 0.3401126408010013
$ws.Range("G6").Value = 0.2802620336086585
$ws.Range("H6").Value = 0.4696667616063799
$ws.Range("I6").Value = 0.1873961009520931
$ws.Range("J6").Value = 0.2995938769134646

$ws.Range("B7").Value = 0.8346807199251247
$ws.Range("C7").Value = 0.03084074308018725
$ws.Range("D7").Value = 0.4693057389143434
$ws.Range("E7").Value = 0.2320642294485963
$ws.Range("F7").Value = 0.1065168539325843
$ws.Range("G7").Value = 0.0667017201231879
$ws.Range("H7").Value = -0.05170824942830832
$ws.Range("I7").Value = 0.0255854095424198
$ws.Range("J7").Value = 0.2036599763872492

$ws.Range("B8").Value = 0.2593551020268173
$ws.Range("C8").Value = 0.3814072148643259
$ws.Range("D8").Value = 0.278620125771039
$ws.Range("E8").Value = 0.04718292575143716
$ws.Range("F8").Value = 0.2792851034485887
$ws.Range("G8").Value = 0.1660068544774337
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.267850820877102
$ws.Range("J8").Value = -0.01918982468696437

$ws.Range("B9").Value = -0.05131399655247695
$ws.Range("C9").Value = -0.1076593250940549
$ws.Range("D9").Value = 0.3021776442138979
$ws.Range("E9").Value = -0.08941583497565135
$ws.Range("F9").Value = 0.9451128337639966
$ws.Range("G9").Value = 0.1246380278546183
$ws.Range("H9").Value = 0.03326093792544518
$ws.Range("I9").Value = -0.103064149378208
$ws.Range("J9").Value = -0.1045867640169459

$ws.Range("B10").Value = 0.3780234968901174
$ws.Range("C10").Value = -0.01964223079621187
$ws.Range("D10").Value = 0.05612065941774814
$ws.Range("E10").Value = 0.09389638537369779
$ws.Range("F10").Value = -0.1065903603684166
$ws.Range("G10").Value = -0.1057360055286801
$ws.Range("H10").Value = 0.09389638537369779
$ws.Range("I10").Value = 0.222427647016995
$ws.Range("J10").Value = -0.1192764273600903

$ws.Range("B11").Value = 0.06784783881846881
$ws.Range("C11").Value = 0.02546097041634578
$ws.Range("D11").Value = 0.05215785028672356
$ws.Range("E11").Value = 0.1551474499143076
$ws.Range("F11").Value = 0.1240487167136861
$ws.Range("G11").Value = 0.0217220622873593
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.03612316378273828
$ws.Range("J11").Value = 0.09626274065685164

$ws.Range("B12").Value = -0.01027960526315787
$ws.Range("C12").Value = -0.03583473861720073
$ws.Range("D12").Value = -0.02574926129168418
$ws.Range("E12").Value = -0.03802281368821287
$ws.Range("F12").Value = -0.02932551319648097
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = -0.03409090909090909
$ws.Range("J12").Value = -0.02704987320371931

$ws.Range("B13").Value = 0.0576923076923077
$ws.Range("C13").Value = -0.01111111111111107
$ws.Range("D13").Value = -0.02162162162162168
$ws.Range("E13").Value = 0.1463414634146342
$ws.Range("F13").Value = -0.03846153846153846
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0.09742120343839548
$ws.Range("J13").Value = -0.06944444444444441

$ws.Range("B14").Value = 0.01282051282051274
$ws.Range("C14").Value = -0.05925925925925932
$ws.Range("D14").Value = 0.04938271604938262
$ws.Range("E14").Value = -0.03314917127071831
$ws.Range("F14").Value = 0.4652777777777778
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.08080110497237562
$ws.Range("J14").Value = 0.06135401974612128

$ws.Range("B15").Value = 0.08803016022620171
$ws.Range("C15").Value = 0.4295938882187375
$ws.Range("D15").Value = 0.2467566390944711
$ws.Range("E15").Value = 0.0007168458781361984
$ws.Range("F15").Value = 0.5412078152753108
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.2053880176919984
$ws.Range("J15").Value = 0.189460555035859

$ws.Range("B16").Value = 0.2296243443743885
$ws.Range("C16").Value = 0.3171130143609565
$ws.Range("D16").Value = 0.4411627132582989
$ws.Range("E16").Value = 0.02701104325362462
$ws.Range("F16").Value = 0.5041558545739103
$ws.Range("G16").Value = 0.05554327312625021
$ws.Range("H16").Value = 0.06607181017159286
$ws.Range("I16").Value = 0.4127498792155826
$ws.Range("J16").Value = 0.2356256661564519

$ws.Range("B17").Value = -0.09018095520617032
$ws.Range("C17").Value = 0.1575185090630584
$ws.Range("D17").Value = 0.1720430107526882
$ws.Range("E17").Value = 0.1380285673879225
$ws.Range("F17").Value = 0.1268575266560322
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0.1449934068400528
$ws.Range("J17").Value = -0.0461798583958576

$ws.Range("B18").Value = 0.2277701310661704
$ws.Range("C18").Value = 0.02891295783279457
$ws.Range("D18").Value = 0.02891295783279457
$ws.Range("E18").Value = 0.07535211267605635
$ws.Range("F18").Value = 0.02891295783279457
$ws.Range("G18").Value = 0.02047608029498984
$ws.Range("H18").Value = 0.01899253161528389
$ws.Range("I18").Value = 0.03845373496969431
$ws.Range("J18").Value = 0.01233779839950207

$ws.Range("B19").Value = 0.4373379463464265
$ws.Range("C19").Value = 0.3957752827261105
$ws.Range("D19").Value = 0.2987437626653427
$ws.Range("E19").Value = 0.1678025151349637
$ws.Range("F19").Value = 0.3643394196220798
$ws.Range("G19").Value = 0.0888019465556701
$ws.Range("H19").Value = 0.02014327193445445
$ws.Range("I19").Value = 0.5053496848550396
$ws.Range("J19").Value = 0.4108260032192536

$ws.Range("B20").Value = 0.3775891341256367
$ws.Range("C20").Value = -0.01464371980676338
$ws.Range("D20").Value = 0.3774012152319131
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0.02558603303608387
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0.000446162998215441
$ws.Range("J20").Value = 0.1557735849056604

$ws.Range("B21").Value = 0.07407713641926475
$ws.Range("C21").Value = 0.5662828008425362
$ws.Range("D21").Value = 0.0289530268830432
$ws.Range("E21").Value = 0.5288960144458497
$ws.Range("F21").Value = -0.006359522004726558
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0.006162187568193721
$ws.Range("J21").Value = -0.01243816254416922

$ws.Range("B22").Value = 0.004969146055828015
$ws.Range("C22").Value = 0.05344262515326211
$ws.Range("D22").Value = 0.01714222394433712
$ws.Range("E22").Value = -0.004671908811878477
$ws.Range("F22").Value = -0.02601210832104633
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0.004190429202395022
$ws.Range("J22").Value = 0.004190429202395022

Write-Output "Updated ARI percentage normed results values (rows 2-22, cols B-J)."
